$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "price" text that looks numeric (e.g. "1.058", "20.527.47",
# "1.472.16"). A leading apostrophe in the literal forces Excel to store the
# value as text (quote-prefixed) -- exactly as if a user typed '<value> into the
# cell -- so trailing zeros / thousand-dot strings are preserved verbatim instead
# of being silently parsed (and rounded/reformatted) as a number.
# In a single-quoted PowerShell string, '' is the escape for a literal single
# quote/apostrophe, so '''20.527.47' is the 11-character text  '20.527.47

$ws.Range("D2").Value = '''20.527.47'
$ws.Range("E2").Value = '  +1.74%  '

$ws.Range("D3").Value = '''1.470.15'
$ws.Range("E3").Value = '  +2.71%  '

$ws.Range("E4").Value = '  +0.64%  '

$ws.Range("D5").Value = '''0.9619'
$ws.Range("E5").Value = '  -3.35%  '

$ws.Range("D6").Value = '''276.27'
$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("D7").Value = '''0.3649'
$ws.Range("E7").Value = '  -1.63%  '

$ws.Range("D8").Value = '''0.3057'
$ws.Range("E8").Value = '  -3.18%  '

$ws.Range("D9").Value = '''40.12'
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("D10").Value = '''1.058'
$ws.Range("E10").Value = '  -0.25%  '

$ws.Range("D11").Value = '''0.06627'
$ws.Range("E11").Value = '  +0.62%  '

$ws.Range("D12").Value = '''1.000'
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = '''5.469'
$ws.Range("E13").Value = '  -1.70%  '

$ws.Range("D14").Value = '''18.12'
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("D15").Value = '''6.176'
$ws.Range("E15").Value = '  -0.91%  '

$ws.Range("D16").Value = '''0.00001029'
$ws.Range("E16").Value = '  -0.34%  '

$ws.Range("D17").Value = '''1.472.74'
$ws.Range("E17").Value = '  +2.75%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '''0.9635'
$ws.Range("E18").Value = '  -3.16%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '''0.05923'
$ws.Range("E19").Value = '  +2.96%  '

$ws.Range("E20").Value = '  -3.96%  '

$ws.Range("D21").Value = '''5.460'
$ws.Range("E21").Value = '  -3.11%  '

$ws.Range("D22").Value = '''14.52'
$ws.Range("E22").Value = '  -2.41%  '

$ws.Range("D23").Value = '''11.06'
$ws.Range("E23").Value = '  -0.86%  '

$ws.Range("D24").Value = '''2.245'
$ws.Range("E24").Value = '  +0.35%  '

$ws.Range("D25").Value = '''20.586.83'
$ws.Range("E25").Value = '  +1.89%  '

$ws.Range("D26").Value = '''140.01'
$ws.Range("E26").Value = '  +3.73%  '

$ws.Range("D27").Value = '''2.130'
$ws.Range("E27").Value = '  -8.14%  '

$ws.Range("D28").Value = '''17.22'
$ws.Range("E28").Value = '  -1.44%  '

$ws.Range("D29").Value = '''1.629.96'
$ws.Range("E29").Value = '  +2.33%  '

$ws.Range("D30").Value = '''113.99'
$ws.Range("E30").Value = '  +1.85%  '

$ws.Range("D31").Value = '''3.952'
$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("D32").Value = '''4.964'
$ws.Range("E32").Value = '  -7.03%  '

$ws.Range("D33").Value = '''0.8113'
$ws.Range("E33").Value = '  -4.22%  '

$ws.Range("D34").Value = '''0.07943'
$ws.Range("E34").Value = '  +1.65%  '

$ws.Range("D35").Value = '''1.538'
$ws.Range("E35").Value = '  +3.12%  '

$ws.Range("E36").Value = '  +8.62%  '

$ws.Range("D37").Value = '''0.05811'
$ws.Range("E37").Value = '  -1.87%  '

$ws.Range("D38").Value = '''4.715'
$ws.Range("E38").Value = '  -4.53%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '''7.653'
$ws.Range("E39").Value = '  -2.77%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.02036'
$ws.Range("E40").Value = '  -1.66%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '''10.44'
$ws.Range("E41").Value = '  -3.44%  '

$ws.Range("D42").Value = '''0.9598'
$ws.Range("E42").Value = '  -3.61%  '

$ws.Range("D43").Value = '''0.1879'
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("D44").Value = '''0.5281'
$ws.Range("E44").Value = '  -1.86%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '''3.509'
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''12.11'
$ws.Range("E46").Value = '  -2.68%  '

$ws.Range("D47").Value = '''117.95'
$ws.Range("E47").Value = '  -1.16%  '

$ws.Range("D48").Value = '''0.5188'
$ws.Range("E48").Value = '  -1.85%  '

$ws.Range("D49").Value = '''1.785'
$ws.Range("E49").Value = '  -1.08%  '

$ws.Range("D50").Value = '''0.06455'
$ws.Range("E50").Value = '  +2.87%  '

$ws.Range("D51").Value = '''0.9893'
$ws.Range("E51").Value = '  -0.82%  '
